$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$bf = New-Object "object[,]" 1,5
$inArr = New-Object "object[,]" 1,6

$bf[0,0] = 1.02
$bf[0,1] = 1.025412041794131
$bf[0,2] = 1.050621060138217
$bf[0,3] = 1.03775829515756
$bf[0,4] = 1.053583496790876
$ws.Range("B2:F2").Value2 = $bf
$inArr[0,0] = 1.03928360183805
$inArr[0,1] = 1.030581209091588
$inArr[0,2] = 1.053374829319516
$inArr[0,3] = 1.040548265998759
$inArr[0,4] = 1.056329064072817
$inArr[0,5] = 1.014205257130815
$ws.Range("I2:N2").Value2 = $inArr

$bf[0,0] = 1.02
$bf[0,1] = 1.026306042787938
$bf[0,2] = 1.05120623247427
$bf[0,3] = 1.038514992205918
$bf[0,4] = 1.054336578937991
$ws.Range("B3:F3").Value2 = $bf
$inArr[0,0] = 1.039384693642305
$inArr[0,1] = 1.031114819076966
$inArr[0,2] = 1.053773070461125
$inArr[0,3] = 1.041114977310231
$inArr[0,4] = 1.056895373408345
$inArr[0,5] = 1.014384928622142
$ws.Range("I3:N3").Value2 = $inArr

$bf[0,0] = 1.02
$bf[0,1] = 1.026885420925439
$bf[0,2] = 1.051584832505319
$bf[0,3] = 1.039005617920432
$bf[0,4] = 1.054824452955891
$ws.Range("B4:F4").Value2 = $bf
$inArr[0,0] = 1.039448788040928
$inArr[0,1] = 1.031460378357532
$inArr[0,2] = 1.054030013300894
$inArr[0,3] = 1.041482005406091
$inArr[0,4] = 1.057261726473567
$inArr[0,5] = 1.014501202929639
$ws.Range("I4:N4").Value2 = $inArr

$bf[0,0] = 1.02
$bf[0,1] = 1.027129205446029
$bf[0,2] = 1.051743982890376
$bf[0,3] = 1.039212112644857
$bf[0,4] = 1.055029691789936
$ws.Range("B5:F5").Value2 = $bf
$inArr[0,0] = 1.039475416809543
$inArr[0,1] = 1.031605716801027
$inArr[0,2] = 1.054137851862167
$inArr[0,3] = 1.041636381088309
$inArr[0,4] = 1.057415718716679
$inArr[0,5] = 1.014550087726847
$ws.Range("I5:N5").Value2 = $inArr

$bf[0,0] = 1.02
$bf[0,1] = 1.027170150444154
$bf[0,2] = 1.051770704111337
$bf[0,3] = 1.039246797784145
$bf[0,4] = 1.055064160241113
$ws.Range("B6:F6").Value2 = $bf
$inArr[0,0] = 1.039479869312917
$inArr[0,1] = 1.031630123560351
$inArr[0,2] = 1.054155947815244
$inArr[0,3] = 1.041662305921965
$inArr[0,4] = 1.057441573315926
$inArr[0,5] = 1.014558295866134
$ws.Range("I6:N6").Value2 = $inArr

$bf[0,0] = 1.02
$bf[0,1] = 1.026888677547249
$bf[0,2] = 1.051586959134623
$bf[0,3] = 1.039008376189328
$bf[0,4] = 1.054827194833203
$ws.Range("B7:F7").Value2 = $bf
$inArr[0,0] = 1.03944914510054
$inArr[0,1] = 1.031462320121784
$inArr[0,2] = 1.054031454955222
$inArr[0,3] = 1.041484067880501
$inArr[0,4] = 1.057263784216437
$inArr[0,5] = 1.014501856119191
$ws.Range("I7:N7").Value2 = $inArr

$bf[0,0] = 1.02
$bf[0,1] = 1.025713986455473
$bf[0,2] = 1.050818830018908
$bf[0,3] = 1.038013817934171
$bf[0,4] = 1.053837883151301
$ws.Range("B8:F8").Value2 = $bf
$inArr[0,0] = 1.039318038861654
$inArr[0,1] = 1.030761486803484
$inArr[0,2] = 1.05350957020516
$inArr[0,3] = 1.040739719957846
$inArr[0,4] = 1.056520467902024
$inArr[0,5] = 1.014265974595212
$ws.Range("I8:N8").Value2 = $inArr

$bf[0,0] = 1.02
$bf[0,1] = 1.02365098531053
$bf[0,2] = 1.049465034084543
$bf[0,3] = 1.036268959183005
$bf[0,4] = 1.052099122458826
$ws.Range("B9:F9").Value2 = $bf
$inArr[0,0] = 1.039076948280336
$inArr[0,1] = 1.02952871147328
$inArr[0,2] = 1.052584303895978
$inArr[0,3] = 1.039430659830828
$inArr[0,4] = 1.055210061963227
$inArr[0,5] = 1.013850455809769
$ws.Range("I9:N9").Value2 = $inArr

$bf[0,0] = 1.02
$bf[0,1] = 1.022280412959304
$bf[0,2] = 1.048562460735517
$bf[0,3] = 1.035110993547989
$bf[0,4] = 1.050943123943906
$ws.Range("B10:F10").Value2 = $bf
$inArr[0,0] = 1.038909502197062
$inArr[0,1] = 1.028708401188176
$inArr[0,2] = 1.051963767496228
$inArr[0,3] = 1.038559771447564
$inArr[0,4] = 1.054336168538048
$inArr[0,5] = 1.013573562025132
$ws.Range("I10:N10").Value2 = $inArr

$bf[0,0] = 1.02
$bf[0,1] = 1.021688087413886
$bf[0,2] = 1.048171653510835
$bf[0,3] = 1.034610855694907
$bf[0,4] = 1.050443343942268
$ws.Range("B11:F11").Value2 = $bf
$inArr[0,0] = 1.038835412620993
$inArr[0,1] = 1.028353578307498
$inArr[0,2] = 1.051694213551874
$inArr[0,3] = 1.038183116690402
$inArr[0,4] = 1.053957716692048
$inArr[0,5] = 1.013453698076597
$ws.Range("I11:N11").Value2 = $inArr

$bf[0,0] = 1.02
$bf[0,1] = 1.021468244241134
$bf[0,2] = 1.048026494457246
$bf[0,3] = 1.03442527475318
$bf[0,4] = 1.050257822016899
$ws.Range("B12:F12").Value2 = $bf
$inArr[0,0] = 1.038807655225668
$inArr[0,1] = 1.028221839141923
$inArr[0,2] = 1.05159396205404
$inArr[0,3] = 1.038043278770138
$inArr[0,4] = 1.05381713708018
$inArr[0,5] = 1.013409180688829
$ws.Range("I12:N12").Value2 = $inArr

$bf[0,0] = 1.02
$bf[0,1] = 1.021515393500385
$bf[0,2] = 1.048057631343948
$bf[0,3] = 1.034465073752692
$bf[0,4] = 1.050297611696771
$ws.Range("B13:F13").Value2 = $bf
$inArr[0,0] = 1.038813620003451
$inArr[0,1] = 1.028250095002544
$inArr[0,2] = 1.051615472067975
$inArr[0,3] = 1.038073271358631
$inArr[0,4] = 1.05384729211526
$inArr[0,5] = 1.013418729566688
$ws.Range("I13:N13").Value2 = $inArr

$bf[0,0] = 1.02
$bf[0,1] = 1.021669911573045
$bf[0,2] = 1.048159654520826
$bf[0,3] = 1.034595511578006
$bf[0,4] = 1.050428006209595
$ws.Range("B14:F14").Value2 = $bf
$inArr[0,0] = 1.038833123020803
$inArr[0,1] = 1.028342687519732
$inArr[0,2] = 1.051685929316123
$inArr[0,3] = 1.038171556239889
$inArr[0,4] = 1.053946096444286
$inArr[0,5] = 1.01345001814211
$ws.Range("I14:N14").Value2 = $inArr

$bf[0,0] = 1.02
$bf[0,1] = 1.021765138212422
$bf[0,2] = 1.048222514982897
$bf[0,3] = 1.034675904195072
$bf[0,4] = 1.050508362356169
$ws.Range("B15:F15").Value2 = $bf
$inArr[0,0] = 1.038845108062321
$inArr[0,1] = 1.028399744528046
$inArr[0,2] = 1.05172932355226
$inArr[0,3] = 1.038232121910307
$inArr[0,4] = 1.054006972344794
$inArr[0,5] = 1.013469296798553
$ws.Range("I15:N15").Value2 = $inArr

$bf[0,0] = 1.02
$bf[0,1] = 1.022319747823927
$bf[0,2] = 1.048588397785455
$bf[0,3] = 1.035144212945634
$bf[0,4] = 1.050976309211439
$ws.Range("B16:F16").Value2 = $bf
$inArr[0,0] = 1.038914385966593
$inArr[0,1] = 1.028731957659538
$inArr[0,2] = 1.051981638979019
$inArr[0,3] = 1.038584778298497
$inArr[0,4] = 1.054361284236181
$inArr[0,5] = 1.013581517734399
$ws.Range("I16:N16").Value2 = $inArr

$bf[0,0] = 1.02
$bf[0,1] = 1.022667946358725
$bf[0,2] = 1.04881791143325
$bf[0,3] = 1.035438311707285
$bf[0,4] = 1.051270048941139
$ws.Range("B17:F17").Value2 = $bf
$inArr[0,0] = 1.038957418562389
$inArr[0,1] = 1.028940447949652
$inArr[0,2] = 1.05213968136798
$inArr[0,3] = 1.038806110688703
$inArr[0,4] = 1.054583522549521
$inArr[0,5] = 1.013651920105498
$ws.Range("I17:N17").Value2 = $inArr

$bf[0,0] = 1.02
$bf[0,1] = 1.022871154697886
$bf[0,2] = 1.048951783967333
$bf[0,3] = 1.035609976844021
$bf[0,4] = 1.051441456983957
$ws.Range("B18:F18").Value2 = $bf
$inArr[0,0] = 1.038982365830371
$inArr[0,1] = 1.029062093090137
$inArr[0,2] = 1.052231781962156
$inArr[0,3] = 1.038935253028296
$inArr[0,4] = 1.054713145431018
$inArr[0,5] = 1.013692987758513
$ws.Range("I18:N18").Value2 = $inArr

$bf[0,0] = 1.02
$bf[0,1] = 1.022940462108666
$bf[0,2] = 1.048997431162566
$bf[0,3] = 1.03566853090703
$bf[0,4] = 1.051499915288617
$ws.Range("B19:F19").Value2 = $bf
$inArr[0,0] = 1.038990846242387
$inArr[0,1] = 1.029103577056099
$inArr[0,2] = 1.052263171779396
$inArr[0,3] = 1.038979294460737
$inArr[0,4] = 1.054757342586616
$inArr[0,5] = 1.013706991285122
$ws.Range("I19:N19").Value2 = $inArr

$bf[0,0] = 1.02
$bf[0,1] = 1.022630576560501
$bf[0,2] = 1.048793286660627
$bf[0,3] = 1.035406745029309
$bf[0,4] = 1.051238525713138
$ws.Range("B20:F20").Value2 = $bf
$inArr[0,0] = 1.03895281738329
$inArr[0,1] = 1.028918075153012
$inArr[0,2] = 1.052122733473532
$inArr[0,3] = 1.038782359369825
$inArr[0,4] = 1.054559678982476
$inArr[0,5] = 1.01364436626901
$ws.Range("I20:N20").Value2 = $inArr

$bf[0,0] = 1.02
$bf[0,1] = 1.021624405067233
$bf[0,2] = 1.048129611105902
$bf[0,3] = 1.034557095559739
$bf[0,4] = 1.050389604983525
$ws.Range("B21:F21").Value2 = $bf
$inArr[0,0] = 1.038827386410542
$inArr[0,1] = 1.028315419724397
$inArr[0,2] = 1.051665184909422
$inArr[0,3] = 1.038142611891306
$inArr[0,4] = 1.053917001176184
$inArr[0,5] = 1.013440804283042
$ws.Range("I21:N21").Value2 = $inArr

$bf[0,0] = 1.02
$bf[0,1] = 1.020992786222813
$bf[0,2] = 1.047712357039628
$bf[0,3] = 1.034024002041748
$bf[0,4] = 1.049856542205952
$ws.Range("B22:F22").Value2 = $bf
$inArr[0,0] = 1.038747150846374
$inArr[0,1] = 1.027936841652401
$inArr[0,2] = 1.051376771082393
$inArr[0,3] = 1.037740773886185
$inArr[0,4] = 1.05351289177802
$inArr[0,5] = 1.013312848376429
$ws.Range("I22:N22").Value2 = $inArr

$bf[0,0] = 1.02
$bf[0,1] = 1.021327523948192
$bf[0,2] = 1.047933548203641
$bf[0,3] = 1.034306498700005
$bf[0,4] = 1.050139062985168
$ws.Range("B23:F23").Value2 = $bf
$inArr[0,0] = 1.038789815033216
$inArr[0,1] = 1.028137500857051
$inArr[0,2] = 1.051529733807085
$inArr[0,3] = 1.037953757728894
$inArr[0,4] = 1.053727120282288
$inArr[0,5] = 1.013380677069051
$ws.Range("I23:N23").Value2 = $inArr

$bf[0,0] = 1.02
$bf[0,1] = 1.022647462023956
$bf[0,2] = 1.048804413532272
$bf[0,3] = 1.035421008274651
$bf[0,4] = 1.051252769472654
$ws.Range("B24:F24").Value2 = $bf
$inArr[0,0] = 1.038954896930652
$inArr[0,1] = 1.028928184344831
$inArr[0,2] = 1.05213039175374
$inArr[0,3] = 1.038793091436085
$inArr[0,4] = 1.054570452879654
$inArr[0,5] = 1.013647779512922
$ws.Range("I24:N24").Value2 = $inArr

$bf[0,0] = 1.02
$bf[0,1] = 1.024183487570641
$bf[0,2] = 1.049815039681719
$bf[0,3] = 1.036719125543692
$bf[0,4] = 1.052548083443307
$ws.Range("B25:F25").Value2 = $bf
$inArr[0,0] = 1.039140463751481
$inArr[0,1] = 1.029847147335111
$inArr[0,2] = 1.052824165420716
$inArr[0,3] = 1.039768769075849
$inArr[0,4] = 1.055548892053561
$inArr[0,5] = 1.013957858514564
$ws.Range("I25:N25").Value2 = $inArr

